$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two outreach rows for Microsoft / IBM (rows 3 and 4), keeping
# only the header row and the Google row. The row/column cells themselves
# stay in place (so the sheet dimension / row structure is unchanged) but
# their content is cleared, leaving just the original style on column A.
$ws.Range("A3:E4").ClearContents()

# The mailto hyperlinks that used to live on A3 (Microsoft) and A4 (IBM)
# must go away as well. This host's Hyperlinks collection only supports
# clearing every hyperlink on the sheet at once, so remove them all and
# then restore the one that should remain (A2 -> hr@google.com), putting
# its cell formatting back the way it was.
$ws.Range("A2:A4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hr@google.com")
$ws.Range("A2").Style = "Hyperlink"

# Update the saved selection to reflect where the user last clicked.
$ws.Range("A3").Select()

Write-Host "Outreach rows for Microsoft and IBM removed; only Google entry remains."
